$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp value already stored in A9
$ws.Range("A9").Value = 44322.77261097686

# Append the new row of data (row 10)
$ws.Range("A10").Value = 44323.77488254345
$ws.Range("B10").Value = 73924
$ws.Range("C10").Value = 62188
$ws.Range("D10").Value = 3194
$ws.Range("E10").Value = 2029
$ws.Range("F10").Value = 1432
$ws.Range("G10").Value = 19303
$ws.Range("H10").Value = 1317
$ws.Range("I10").Value = 820
$ws.Range("J10").Value = 205

# New date cell should carry the same number format (s="2") as the rest of column A
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat
